$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Correction": values previously in columns I (Avg_Throughput_Before) and
# J (Avg_Throughput_After) for data rows 2-22 were off by a factor of 10 -
# multiply them back up.
for ($r = 2; $r -le 22; $r++) {
    $ws.Cells.Item($r, 9).Value2 = $ws.Cells.Item($r, 9).Value2 * 10
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($r, 10).Value2 * 10
}

# Column L picked up an explicit (best-fit) width of 10 characters.
$ws.Columns.Item(12).ColumnWidth = 9.1666666666667

# Active selection moved to O6.
$ws.Range("O6").Select() | Out-Null
